{"js": "// Update the date line and the 20x5 table of arithmetic answers.\nconst body = context.document.body;\n\n// 1) Update the date paragraph (the first paragraph in the body, before the table).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\n\nif (dateParagraph.text === \"2025-04-07 Monday\") {\n  dateParagraph.insertText(\"2025-04-08 Tuesday\", Word.InsertLocation.replace);\n}\n\n// 2) Update every cell of the table with its new arithmetic expression, in row-major\n//    order, matching the table's existing 20x5 layout. Assigning `table.values`\n//    rewrites only the text of each cell's first run and preserves formatting\n//    (fonts, size, paragraph alignment, etc.).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = [\n  [\"80-26=54\", \"40+27=67\", \"24+60=84\", \"31+53=84\", \"92-47=45\"],\n  [\"56+19=75\", \"36+51=87\", \"63-53=10\", \"67+2=69\", \"19+33=52\"],\n  [\"0+88=88\", \"12+85=97\", \"10+64=74\", \"57-55=2\", \"72-58=14\"],\n  [\"45-19=26\", \"58-57=1\", \"1+8=9\", \"78-23=55\", \"92-84=8\"],\n  [\"41+27=68\", \"66-1=65\", \"63+34=97\", \"30+44=74\", \"66-26=40\"],\n  [\"29+10=39\", \"11+4=15\", \"10+53=63\", \"60+9=69\", \"64+15=79\"],\n  [\"96-36=60\", \"75-62=13\", \"14+59=73\", \"24+20=44\", \"21+49=70\"],\n  [\"48+31=79\", \"80-14=66\", \"12+51=63\", \"36+0=36\", \"23+0=23\"],\n  [\"29+18=47\", \"7+33=40\", \"81-11=70\", \"99-57=42\", \"55+35=90\"],\n  [\"67+32=99\", \"8+67=75\", \"64-31=33\", \"30-15=15\", \"29+29=58\"],\n  [\"5+94=99\", \"5+43=48\", \"44+47=91\", \"58-34=24\", \"69-53=16\"],\n  [\"82-1=81\", \"31-27=4\", \"78-35=43\", \"25+12=37\", \"55-44=11\"],\n  [\"86-33=53\", \"54-41=13\", \"42+33=75\", \"30-5=25\", \"17+59=76\"],\n  [\"68-56=12\", \"77-4=73\", \"95-77=18\", \"49-2=47\", \"86-26=60\"],\n  [\"64-21=43\", \"55+10=65\", \"8+35=43\", \"32+25=57\", \"58+33=91\"],\n  [\"50+27=77\", \"5+13=18\", \"53+11=64\", \"8+27=35\", \"85-23=62\"],\n  [\"10+75=85\", \"37+6=43\", \"58-16=42\", \"28+8=36\", \"49+3=52\"],\n  [\"86-83=3\", \"88-51=37\", \"76-30=46\", \"20+74=94\", \"2+96=98\"],\n  [\"39-33=6\", \"88-17=71\", \"41+1=42\", \"98-29=69\", \"56-48=8\"],\n  [\"97-37=60\", \"96-59=37\", \"19+76=95\", \"31+60=91\", \"63+26=89\"]\n];\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the date line (first paragraph, before the table).\n$dateParagraph = $d.Paragraphs.Item(1)\nif ($dateParagraph.Range.Text.StartsWith(\"2025-04-07 Monday\")) {\n    $dateParagraph.Range.Text = \"2025-04-08 Tuesday\"\n}\n\n# 2) Update every cell of the 20x5 table with the new arithmetic expressions,\n#    preserving existing cell formatting (Range.Text only rewrites the text\n#    before the end-of-cell mark).\n$t = $d.Tables.Item(1)\n\n$values = @(\n    @(\"80-26=54\", \"40+27=67\", \"24+60=84\", \"31+53=84\", \"92-47=45\"),\n    @(\"56+19=75\", \"36+51=87\", \"63-53=10\", \"67+2=69\", \"19+33=52\"),\n    @(\"0+88=88\", \"12+85=97\", \"10+64=74\", \"57-55=2\", \"72-58=14\"),\n    @(\"45-19=26\", \"58-57=1\", \"1+8=9\", \"78-23=55\", \"92-84=8\"),\n    @(\"41+27=68\", \"66-1=65\", \"63+34=97\", \"30+44=74\", \"66-26=40\"),\n    @(\"29+10=39\", \"11+4=15\", \"10+53=63\", \"60+9=69\", \"64+15=79\"),\n    @(\"96-36=60\", \"75-62=13\", \"14+59=73\", \"24+20=44\", \"21+49=70\"),\n    @(\"48+31=79\", \"80-14=66\", \"12+51=63\", \"36+0=36\", \"23+0=23\"),\n    @(\"29+18=47\", \"7+33=40\", \"81-11=70\", \"99-57=42\", \"55+35=90\"),\n    @(\"67+32=99\", \"8+67=75\", \"64-31=33\", \"30-15=15\", \"29+29=58\"),\n    @(\"5+94=99\", \"5+43=48\", \"44+47=91\", \"58-34=24\", \"69-53=16\"),\n    @(\"82-1=81\", \"31-27=4\", \"78-35=43\", \"25+12=37\", \"55-44=11\"),\n    @(\"86-33=53\", \"54-41=13\", \"42+33=75\", \"30-5=25\", \"17+59=76\"),\n    @(\"68-56=12\", \"77-4=73\", \"95-77=18\", \"49-2=47\", \"86-26=60\"),\n    @(\"64-21=43\", \"55+10=65\", \"8+35=43\", \"32+25=57\", \"58+33=91\"),\n    @(\"50+27=77\", \"5+13=18\", \"53+11=64\", \"8+27=35\", \"85-23=62\"),\n    @(\"10+75=85\", \"37+6=43\", \"58-16=42\", \"28+8=36\", \"49+3=52\"),\n    @(\"86-83=3\", \"88-51=37\", \"76-30=46\", \"20+74=94\", \"2+96=98\"),\n    @(\"39-33=6\", \"88-17=71\", \"41+1=42\", \"98-29=69\", \"56-48=8\"),\n    @(\"97-37=60\", \"96-59=37\", \"19+76=95\", \"31+60=91\", \"63+26=89\"),\n)\n\nfor ($r = 1; $r -le $values.Count; $r++) {\n    $row = $values[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n"}
